# Apply cryptos list update (prices / volume changes + 4 row reorderings)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.199.44"
$ws.Range("E2").Value = "  -13.97%  "

$ws.Range("D3").Value = "2.311.42"
$ws.Range("E3").Value = "  -20.45%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "'434.26"
$ws.Range("E5").Value = "  -17.80%  "

$ws.Range("D6").Value = "'121.43"
$ws.Range("E6").Value = "  -15.19%  "

$ws.Range("D7").Value = "'0.995"
$ws.Range("E7").Value = "  -0.45%  "

$ws.Range("D8").Value = "'0.464"
$ws.Range("E8").Value = "  -16.24%  "

$ws.Range("D9").Value = "2.305.58"
$ws.Range("E9").Value = "  -20.85%  "

$ws.Range("D10").Value = "'5.17"
$ws.Range("E10").Value = "  -14.11%  "

$ws.Range("D11").Value = "'0.0874"
$ws.Range("E11").Value = "  -19.14%  "

$ws.Range("D12").Value = "'0.300"
$ws.Range("E12").Value = "  -16.68%  "

$ws.Range("D13").Value = "'0.121"
$ws.Range("E13").Value = "  -5.53%  "

$ws.Range("D14").Value = "52.201.24"
$ws.Range("E14").Value = "  -13.93%  "

$ws.Range("D15").Value = "'18.75"
$ws.Range("E15").Value = "  -17.86%  "

$ws.Range("D16").Value = "'0.0000118"
$ws.Range("E16").Value = "  -16.60%  "

$ws.Range("D17").Value = "2.318.08"
$ws.Range("E17").Value = "  -20.43%  "

$ws.Range("D18").Value = "'3.92"
$ws.Range("E18").Value = "  -21.72%  "

$ws.Range("D19").Value = "'294.02"
$ws.Range("E19").Value = "  -18.37%  "

$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'8.85"
$ws.Range("E20").Value = "  -24.32%  "

$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'0.998"
$ws.Range("E21").Value = "  -0.13%  "

$ws.Range("D22").Value = "'5.61"
$ws.Range("E22").Value = "  -1.33%  "

$ws.Range("D23").Value = "'5.17"
$ws.Range("E23").Value = "  -22.48%  "

$ws.Range("D24").Value = "'53.03"
$ws.Range("E24").Value = "  -18.25%  "

$ws.Range("D25").Value = "'0.362"
$ws.Range("E25").Value = "  -20.42%  "

$ws.Range("D26").Value = "'0.144"
$ws.Range("E26").Value = "  -20.45%  "

$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "'6.85"
$ws.Range("E27").Value = "  -13.31%  "

$ws.Range("B28").Value = "USDe"
$ws.Range("C28").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D28").Value = "'0.994"
$ws.Range("E28").Value = "  -0.48%  "

$ws.Range("D29").Value = "0.0₃0669"
$ws.Range("E29").Value = "  -21.23%  "

$ws.Range("D30").Value = "'141.82"
$ws.Range("E30").Value = "  -6.09%  "

$ws.Range("D31").Value = "'16.96"
$ws.Range("E31").Value = "  -14.34%  "

$ws.Range("D32").Value = "'1.32"
$ws.Range("E32").Value = "  -21.71%  "

$ws.Range("D33").Value = "'4.67"
$ws.Range("E33").Value = "  -16.44%  "

$ws.Range("D34").Value = "'3.45"
$ws.Range("E34").Value = "  -20.86%  "

$ws.Range("D35").Value = "'0.809"
$ws.Range("E35").Value = "  -19.54%  "

$ws.Range("E36").Value = "  -0.34%  "

$ws.Range("D37").Value = "'0.990"
$ws.Range("E37").Value = "  -17.66%  "

$ws.Range("D38").Value = "'31.74"
$ws.Range("E38").Value = "  -16.25%  "

$ws.Range("D39").Value = "'10.15"
$ws.Range("E39").Value = "  -1.66%  "

$ws.Range("D40").Value = "'3.12"
$ws.Range("E40").Value = "  -16.06%  "

$ws.Range("D41").Value = "'1.20"
$ws.Range("E41").Value = "  -18.91%  "

$ws.Range("D42").Value = "'0.0497"
$ws.Range("E42").Value = "  -15.09%  "

$ws.Range("D43").Value = "1.881.19"
$ws.Range("E43").Value = "  -18.12%  "

$ws.Range("D44").Value = "'0.522"
$ws.Range("E44").Value = "  -19.74%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0205"
$ws.Range("E45").Value = "  -14.23%  "

$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.0818"
$ws.Range("E46").Value = "  -11.08%  "

$ws.Range("D47").Value = "'15.56"
$ws.Range("E47").Value = "  -23.85%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'3.88"
$ws.Range("E48").Value = "  -21.42%  "

$ws.Range("B49").Value = "ZEEBU"
$ws.Range("C49").Value = "https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu"
$ws.Range("D49").Value = "'4.55"
$ws.Range("E49").Value = "  -5.88%  "

$ws.Range("D50").Value = "'14.95"
$ws.Range("E50").Value = "  -18.56%  "

$ws.Range("D51").Value = "'4.38"
$ws.Range("E51").Value = "  -15.52%  "

